# Updated output for 2020 run reconstruction
$wb = $excel.ActiveWorkbook

# --- Sheet: Escapement ---
# Only a new row 59 is appended; rows 57/58 are unchanged.
$ws = $wb.Worksheets.Item("Escapement")
$ws.Range("A59").Value = 323809.204451901
$ws.Range("B59").Value = 2243895.849199
$ws.Range("C59").Value = 1228072.698399

# --- Sheet: Total Catch ---
# Rows 57 and 58 updated, row 59 appended.
$ws = $wb.Worksheets.Item("Total Catch")
$ws.Range("A57").Value = 1149316.32227137
$ws.Range("B57").Value = 15172918.7120943
$ws.Range("C57").Value = 8354718.93379377

$ws.Range("A58").Value = 1102852.59372615
$ws.Range("B58").Value = 10267464.1105131
$ws.Range("C58").Value = 3597165.06696538

$ws.Range("A59").Value = 898133.45050504
$ws.Range("B59").Value = 5851856.86129701
$ws.Range("C59").Value = 2442312.033493

# --- Sheet: Run Size ---
# Rows 57 and 58 updated, row 59 appended.
$ws = $wb.Worksheets.Item("Run Size")
$ws.Range("A57").Value = 1920088.32231608
$ws.Range("B57").Value = 22680172.7116691
$ws.Range("C57").Value = 9602178.93404677

$ws.Range("A58").Value = 1358926.59367825
$ws.Range("B58").Value = 12340740.1103681
$ws.Range("C58").Value = 4306596.06695639

$ws.Range("A59").Value = 1221947.45055694
$ws.Range("B59").Value = 8095742.86049601
$ws.Range("C59").Value = 3670371.033892

# --- Sheet: Run Size no Offshore ---
# Row 58 updated, row 59 appended.
$ws = $wb.Worksheets.Item("Run Size no Offshore")
$ws.Range("A58").Value = 1347502.3239042
$ws.Range("B58").Value = 12198050.7831206
$ws.Range("C58").Value = 4255997.519754

$ws.Range("A59").Value = 1205082.9590016
$ws.Range("B59").Value = 7950445.03944901
$ws.Range("C59").Value = 3602937.819464
